# Trade #87 closed at 2026-02-18 00:29:11 - unknown UNKNOWN +0.000%
#
# This script applies the commit's changes to live_trading_results.xlsx:
#   - Summary sheet: refresh aggregate metrics
#   - Strategy Status sheet: refresh MarketMaking strategy row
#   - All Trades sheet: close out trade #115 (MarketMaking) early, and append
#     two brand-new OPEN trades (#144 HighProbConvergence, #145 MarketMaking)
#   - HighProbConvergence sheet: append the strategy's own copy of trade #144
#   - MarketMaking sheet: close out its copy of trade #115, append trade #145

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Writes a literal text value into a cell without Excel's automatic
    # number/date/time detection turning it into a serial number, and
    # without leaving behind a lingering text number-format on the cell.
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-BlankCell {
    # Recreates an "empty" cell (mirrors the original inlineStr-with-no-text
    # placeholder cells used for not-yet-populated columns like Exit Price).
    param($cell)
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.42
$summary.Range("B4").Value = 0.53
$summary.Range("B6").Value = 115
$summary.Range("B7").Value = 55
$summary.Range("B9").Value = 47.83

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.5
$status.Range("D6").Value = 35
$status.Range("E6").Value = -0.31
$status.Range("F6").Value = -0.5
$status.Range("G6").Value = 48.57

# ---------------------------------------------------------------------
# All Trades sheet - close trade #115 (MarketMaking) early, row 116
# Columns: A#, B Date, C Time, D Strategy, E Side, F EntryPrice,
#          G ExitPrice, H Status, I P&L%, J P&L$, K CapitalAfter,
#          L ExitReason, M Duration, N EntrySlip, O ExitSlip,
#          P Confidence, Q EntryReason
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(116, 7).Value = 0.98
Set-TextCell $allTrades.Cells.Item(116, 8) "CLOSED"
$allTrades.Cells.Item(116, 9).Value = 2.0833
$allTrades.Cells.Item(116, 10).Value = 0.02
$allTrades.Cells.Item(116, 11).Value = 99.5
Set-TextCell $allTrades.Cells.Item(116, 12) "early_exit"
$allTrades.Cells.Item(116, 13).Value = 0.17

# New row 145 - trade #144, HighProbConvergence, still OPEN
$allTrades.Cells.Item(145, 1).Value = 144
Set-TextCell $allTrades.Cells.Item(145, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(145, 3) "00:29:05"
Set-TextCell $allTrades.Cells.Item(145, 4) "HighProbConvergence"
Set-TextCell $allTrades.Cells.Item(145, 5) "DOWN"
$allTrades.Cells.Item(145, 6).Value = 0.96
Set-BlankCell $allTrades.Cells.Item(145, 7)
Set-TextCell $allTrades.Cells.Item(145, 8) "OPEN"
$allTrades.Cells.Item(145, 9).Value = 0
$allTrades.Cells.Item(145, 10).Value = 0
$allTrades.Cells.Item(145, 11).Value = 100.4130057263667
Set-BlankCell $allTrades.Cells.Item(145, 12)
$allTrades.Cells.Item(145, 13).Value = 0
$allTrades.Cells.Item(145, 14).Value = 0
$allTrades.Cells.Item(145, 15).Value = 0
$allTrades.Cells.Item(145, 16).Value = 0.95
Set-TextCell $allTrades.Cells.Item(145, 17) "Mean reversion DOWN: price 2.19% above mean (z=1.73)"

# New row 146 - trade #145, MarketMaking, still OPEN
$allTrades.Cells.Item(146, 1).Value = 145
Set-TextCell $allTrades.Cells.Item(146, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(146, 3) "00:29:05"
Set-TextCell $allTrades.Cells.Item(146, 4) "MarketMaking"
Set-TextCell $allTrades.Cells.Item(146, 5) "DOWN"
$allTrades.Cells.Item(146, 6).Value = 0.96
Set-BlankCell $allTrades.Cells.Item(146, 7)
Set-TextCell $allTrades.Cells.Item(146, 8) "OPEN"
$allTrades.Cells.Item(146, 9).Value = 0
$allTrades.Cells.Item(146, 10).Value = 0
$allTrades.Cells.Item(146, 11).Value = 99.47967800952271
Set-BlankCell $allTrades.Cells.Item(146, 12)
$allTrades.Cells.Item(146, 13).Value = 0
$allTrades.Cells.Item(146, 14).Value = 0
$allTrades.Cells.Item(146, 15).Value = 0
$allTrades.Cells.Item(146, 16).Value = 0.6
Set-TextCell $allTrades.Cells.Item(146, 17) "Normal spread capture: 190 bps"

# ---------------------------------------------------------------------
# HighProbConvergence sheet - append trade #144 as row 18
# Columns: A#, B Date, C Time, D Strategy, E Side, F EntryPrice,
#          G ExitPrice, H Status, I P&L%, J P&L$, K CapitalAfter,
#          L EntrySlip, M ExitSlip, N Confidence, O EntryReason,
#          P ExitReason, Q Duration
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")

$hpc.Cells.Item(18, 1).Value = 144
Set-TextCell $hpc.Cells.Item(18, 2) "2026-02-18"
Set-TextCell $hpc.Cells.Item(18, 3) "00:29:05"
Set-TextCell $hpc.Cells.Item(18, 4) "HighProbConvergence"
Set-TextCell $hpc.Cells.Item(18, 5) "DOWN"
$hpc.Cells.Item(18, 6).Value = 0.96
Set-BlankCell $hpc.Cells.Item(18, 7)
Set-TextCell $hpc.Cells.Item(18, 8) "OPEN"
$hpc.Cells.Item(18, 9).Value = 0
$hpc.Cells.Item(18, 10).Value = 0
$hpc.Cells.Item(18, 11).Value = 100.4130057263667
$hpc.Cells.Item(18, 12).Value = 0
$hpc.Cells.Item(18, 13).Value = 0
$hpc.Cells.Item(18, 14).Value = 0.95
Set-TextCell $hpc.Cells.Item(18, 15) "Mean reversion DOWN: price 2.19% above mean (z=1.73)"
Set-BlankCell $hpc.Cells.Item(18, 16)
$hpc.Cells.Item(18, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close trade #115 early (row 36), append #145 (row 60)
# Columns: A#, B Date, C Time, D Strategy, E Side, F EntryPrice,
#          G ExitPrice, H Status, I P&L%, J P&L$, K CapitalAfter,
#          L EntrySlip, M ExitSlip, N Confidence, O EntryReason,
#          P ExitReason, Q Duration
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(36, 7).Value = 0.98
Set-TextCell $mm.Cells.Item(36, 8) "CLOSED"
$mm.Cells.Item(36, 9).Value = 2.0833
$mm.Cells.Item(36, 10).Value = 0.02
$mm.Cells.Item(36, 11).Value = 99.5
Set-TextCell $mm.Cells.Item(36, 16) "early_exit"
$mm.Cells.Item(36, 17).Value = 0.17

# New row 60 - trade #145, MarketMaking, still OPEN
$mm.Cells.Item(60, 1).Value = 145
Set-TextCell $mm.Cells.Item(60, 2) "2026-02-18"
Set-TextCell $mm.Cells.Item(60, 3) "00:29:05"
Set-TextCell $mm.Cells.Item(60, 4) "MarketMaking"
Set-TextCell $mm.Cells.Item(60, 5) "DOWN"
$mm.Cells.Item(60, 6).Value = 0.96
Set-BlankCell $mm.Cells.Item(60, 7)
Set-TextCell $mm.Cells.Item(60, 8) "OPEN"
$mm.Cells.Item(60, 9).Value = 0
$mm.Cells.Item(60, 10).Value = 0
$mm.Cells.Item(60, 11).Value = 99.47967800952271
$mm.Cells.Item(60, 12).Value = 0
$mm.Cells.Item(60, 13).Value = 0
$mm.Cells.Item(60, 14).Value = 0.6
Set-TextCell $mm.Cells.Item(60, 15) "Normal spread capture: 190 bps"
Set-BlankCell $mm.Cells.Item(60, 16)
$mm.Cells.Item(60, 17).Value = 0

Write-Output "Edit complete"
